$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1535087719298246
$ws.Range("C2").Value = 0.5921052631578947
$ws.Range("J2").Value = 0.008771929824561403
$ws.Range("P2").Value = 0.1140350877192982
$ws.Range("S2").Value = 0.131578947368421
$ws.Range("B3").Value = 0.04137931034482759
$ws.Range("C3").Value = 0.06206896551724138
$ws.Range("J3").Value = 0.01379310344827586
$ws.Range("P3").Value = 0.7517241379310344
$ws.Range("S3").Value = 0.1310344827586207
$ws.Range("J4").Value = 0.06451612903225806
$ws.Range("P4").Value = 0.6774193548387096
$ws.Range("S4").Value = 0.2580645161290323
$ws.Range("B6").Value = 0.03864734299516908
$ws.Range("D6").Value = 0.004830917874396135
$ws.Range("E6").Value = 0.004830917874396135
$ws.Range("F6").Value = 0.03381642512077294
$ws.Range("J6").Value = 0.2801932367149759
$ws.Range("O6").Value = 0.02898550724637681
$ws.Range("Q6").Value = 0.1739130434782609
$ws.Range("R6").Value = 0.1014492753623188
$ws.Range("S6").Value = 0.3333333333333333
$ws.Range("B7").Value = 0.111731843575419
$ws.Range("D7").Value = 0.00558659217877095
$ws.Range("E7").Value = 0.0111731843575419
$ws.Range("F7").Value = 0.03910614525139665
$ws.Range("J7").Value = 0.2122905027932961
$ws.Range("O7").Value = 0.05027932960893855
$ws.Range("Q7").Value = 0.1229050279329609
$ws.Range("R7").Value = 0.106145251396648
$ws.Range("S7").Value = 0.3407821229050279
$ws.Range("B8").Value = 0.08240534521158129
$ws.Range("D8").Value = 0.0111358574610245
$ws.Range("E8").Value = 0.0022271714922049
$ws.Range("F8").Value = 0.0467706013363029
$ws.Range("J8").Value = 0.1380846325167038
$ws.Range("O8").Value = 0.0133630289532294
$ws.Range("Q8").Value = 0.1915367483296214
$ws.Range("R8").Value = 0.1269487750556793
$ws.Range("S8").Value = 0.3875278396436526
$ws.Range("B9").Value = 0.1180124223602484
$ws.Range("D9").Value = 0.006211180124223602
$ws.Range("F9").Value = 0.07453416149068323
$ws.Range("J9").Value = 0.1118012422360248
$ws.Range("O9").Value = 0.03105590062111801
$ws.Range("Q9").Value = 0.1925465838509317
$ws.Range("R9").Value = 0.1180124223602484
$ws.Range("S9").Value = 0.3478260869565217
$ws.Range("B10").Value = 0.07819548872180451
$ws.Range("D10").Value = 0.01729323308270677
$ws.Range("E10").Value = 0.002255639097744361
$ws.Range("F10").Value = 0.07142857142857142
$ws.Range("J10").Value = 0.1390977443609022
$ws.Range("O10").Value = 0.01353383458646617
$ws.Range("Q10").Value = 0.2045112781954887
$ws.Range("R10").Value = 0.1330827067669173
$ws.Range("S10").Value = 0.3406015037593985
$ws.Range("G11").Value = 0.1160409556313993
$ws.Range("J11").Value = 0.1228668941979522
$ws.Range("K11").Value = 0.204778156996587
$ws.Range("L11").Value = 0.5221843003412969
$ws.Range("S11").Value = 0.03412969283276451
$ws.Range("G12").Value = 0.7349397590361446
$ws.Range("J12").Value = 0.1686746987951807
$ws.Range("K12").Value = 0.01204819277108434
$ws.Range("L12").Value = 0.03614457831325301
$ws.Range("S12").Value = 0.04819277108433735
$ws.Range("G13").Value = 0.7105263157894737
$ws.Range("J13").Value = 0.1842105263157895
$ws.Range("S13").Value = 0.1052631578947368
$ws.Range("F15").Value = 0.03347280334728033
$ws.Range("H15").Value = 0.1548117154811715
$ws.Range("I15").Value = 0.04602510460251046
$ws.Range("J15").Value = 0.393305439330544
$ws.Range("K15").Value = 0.07112970711297072
$ws.Range("M15").Value = 0.008368200836820083
$ws.Range("N15").Value = 0.004184100418410041
$ws.Range("O15").Value = 0.03765690376569038
$ws.Range("S15").Value = 0.2510460251046025
$ws.Range("H16").Value = 0.1623376623376623
$ws.Range("I16").Value = 0.06493506493506493
$ws.Range("J16").Value = 0.4545454545454545
$ws.Range("K16").Value = 0.06493506493506493
$ws.Range("M16").Value = 0.02597402597402598
$ws.Range("O16").Value = 0.05194805194805195
$ws.Range("S16").Value = 0.1753246753246753
$ws.Range("F17").Value = 0.022271714922049
$ws.Range("H17").Value = 0.1826280623608018
$ws.Range("I17").Value = 0.08240534521158129
$ws.Range("J17").Value = 0.4432071269487751
$ws.Range("K17").Value = 0.09354120267260579
$ws.Range("M17").Value = 0.0200445434298441
$ws.Range("O17").Value = 0.0623608017817372
$ws.Range("S17").Value = 0.09354120267260579
$ws.Range("F18").Value = 0.0136518771331058
$ws.Range("H18").Value = 0.2081911262798635
$ws.Range("I18").Value = 0.08532423208191127
$ws.Range("J18").Value = 0.447098976109215
$ws.Range("K18").Value = 0.05460750853242321
$ws.Range("M18").Value = 0.01706484641638225
$ws.Range("N18").Value = 0.003412969283276451
$ws.Range("O18").Value = 0.07508532423208192
$ws.Range("S18").Value = 0.09556313993174062
$ws.Range("F19").Value = 0.01730103806228374
$ws.Range("H19").Value = 0.2119377162629758
$ws.Range("I19").Value = 0.06660899653979238
$ws.Range("J19").Value = 0.3503460207612457
$ws.Range("K19").Value = 0.1280276816608996
$ws.Range("M19").Value = 0.01730103806228374
$ws.Range("O19").Value = 0.08564013840830449
$ws.Range("S19").Value = 0.1228373702422145
